$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 298.1
$ws.Range("J9").Value = 369.18182
$ws.Range("L9").Value = 369.18182
$ws.Range("N9").Value = -707.18182
$ws.Range("H55").Value = 580.55
$ws.Range("I55").Value = 622.2308
$ws.Range("J55").Value = 503.14285
$ws.Range("K55").Value = 622.2308
$ws.Range("L55").Value = 503.14285
$ws.Range("M55").Value = -408.2308
$ws.Range("N55").Value = -931.14285
$ws.Range("H101").Value = 650
$ws.Range("J101").Value = 622
$ws.Range("L101").Value = 1866
$ws.Range("N101").Value = -5110
$ws.Range("H106").Value = 2656.2917
$ws.Range("I106").Value = 2878.3125
$ws.Range("J106").Value = 2212.25
$ws.Range("K106").Value = 2878.3125
$ws.Range("L106").Value = 2212.25
$ws.Range("M106").Value = -2247.3125
$ws.Range("N106").Value = -3474.25
$ws.Range("H107").Value = 774.2727
$ws.Range("I107").Value = 768.2
$ws.Range("K107").Value = 768.2
$ws.Range("M107").Value = 1151.8
$ws.Range("H132").Value = 13846.238
$ws.Range("I132").Value = 1584.919
$ws.Range("J132").Value = 104580
$ws.Range("K132").Value = 4754.757000000001
$ws.Range("L132").Value = 313740
$ws.Range("M132").Value = -2224.757000000001
$ws.Range("N132").Value = -318800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 9499.5
$ws.Range("I21").Value = 999
$ws.Range("J21").Value = 18000
$ws.Range("K21").Value = 999
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = -625
$ws.Range("N21").Value = -18748
$ws.Range("H32").Value = 3479.6072
$ws.Range("I32").Value = 4133.143
$ws.Range("J32").Value = 1519
$ws.Range("K32").Value = 4133.143
$ws.Range("L32").Value = 1519
$ws.Range("M32").Value = -3846.143
$ws.Range("N32").Value = -2093

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H105").Value = 2086.3157
$ws.Range("I105").Value = 1171.909
$ws.Range("K105").Value = 1171.909
$ws.Range("M105").Value = 575.0909999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2886.36
$ws.Range("I31").Value = 1978.4286
$ws.Range("J31").Value = 4041.9092
$ws.Range("K31").Value = 1978.4286
$ws.Range("L31").Value = 4041.9092
$ws.Range("M31").Value = -1683.4286
$ws.Range("N31").Value = -4631.9092
$ws.Range("H34").Value = 2886.36
$ws.Range("I34").Value = 1978.4286
$ws.Range("J34").Value = 4041.9092
$ws.Range("K34").Value = 1978.4286
$ws.Range("L34").Value = 4041.9092
$ws.Range("M34").Value = -1776.4286
$ws.Range("N34").Value = -4445.9092
$ws.Range("H96").Value = 9878.5
$ws.Range("J96").Value = 9878.5
$ws.Range("L96").Value = 9878.5
$ws.Range("N96").Value = -15370.5
$ws.Range("H134").Value = 3476.9666
$ws.Range("I134").Value = 3929.5833
$ws.Range("J134").Value = 1666.5
$ws.Range("K134").Value = 11788.7499
$ws.Range("L134").Value = 4999.5
$ws.Range("M134").Value = -9253.749899999999
$ws.Range("N134").Value = -10069.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 344.7143
$ws.Range("I14").Value = 344.7143
$ws.Range("K14").Value = 1034.1429
$ws.Range("M14").Value = -861.1428999999998
$ws.Range("H25").Value = 1206.6666
$ws.Range("I25").Value = 871.4286
$ws.Range("J25").Value = 1500
$ws.Range("K25").Value = 2614.2858
$ws.Range("L25").Value = 4500
$ws.Range("M25").Value = -2445.2858
$ws.Range("N25").Value = -4838
$ws.Range("H29").Value = 7778280.5
$ws.Range("I29").Value = 12963550
$ws.Range("K29").Value = 38890650
$ws.Range("M29").Value = -38890373
$ws.Range("H30").Value = 1206.6666
$ws.Range("I30").Value = 871.4286
$ws.Range("J30").Value = 1500
$ws.Range("K30").Value = 2614.2858
$ws.Range("L30").Value = 4500
$ws.Range("M30").Value = -2512.2858
$ws.Range("N30").Value = -4704
$ws.Range("H94").Value = 4131.6665
$ws.Range("I94").Value = 2395
$ws.Range("K94").Value = 7185
$ws.Range("M94").Value = -6509
$ws.Range("H107").Value = 1942.1666
$ws.Range("J107").Value = 1827.6923
$ws.Range("L107").Value = 5483.0769
$ws.Range("N107").Value = -9323.0769
$ws.Range("H108").Value = 1361.6666
$ws.Range("I108").Value = 1361.6666
$ws.Range("K108").Value = 4084.9998
$ws.Range("M108").Value = -1204.9998
$ws.Range("H122").Value = 4391.95
$ws.Range("J122").Value = 5738.6665
$ws.Range("L122").Value = 51647.9985
$ws.Range("N122").Value = -56547.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 29909.092
$ws.Range("J88").Value = 29909.092
$ws.Range("L88").Value = 29909.092
$ws.Range("N88").Value = -30811.092
$ws.Range("H91").Value = 29909.092
$ws.Range("J91").Value = 29909.092
$ws.Range("L91").Value = 29909.092
$ws.Range("N91").Value = -33029.092
$ws.Range("H102").Value = 5127.6816
$ws.Range("I102").Value = 4590.95
$ws.Range("K102").Value = 4590.95
$ws.Range("M102").Value = -2968.95
$ws.Range("H104").Value = 44992.5
$ws.Range("J104").Value = 44992.5
$ws.Range("L104").Value = 44992.5
$ws.Range("N104").Value = -51980.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3192.4666
$ws.Range("I46").Value = 3124.75
$ws.Range("J46").Value = 3217.0908
$ws.Range("K46").Value = 3124.75
$ws.Range("L46").Value = 3217.0908
$ws.Range("M46").Value = -2936.75
$ws.Range("N46").Value = -3593.0908
$ws.Range("H94").Value = 55000
$ws.Range("J94").Value = 55000
$ws.Range("L94").Value = 55000
$ws.Range("N94").Value = -56352
$ws.Range("H132").Value = 2184.6155
$ws.Range("I132").Value = 1712.8
$ws.Range("J132").Value = 3757.3333
$ws.Range("K132").Value = 5138.4
$ws.Range("L132").Value = 11271.9999
$ws.Range("M132").Value = -2608.4
$ws.Range("N132").Value = -16331.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3450
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 4083.3333
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 4083.3333
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -5331.3333
$ws.Range("H65").Value = 3450
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 4083.3333
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 20416.6665
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -26656.6665
$ws.Range("H95").Value = 35743.8
$ws.Range("J95").Value = 35743.8
$ws.Range("L95").Value = 35743.8
$ws.Range("N95").Value = -41235.8
$ws.Range("H132").Value = 5011.614
$ws.Range("I132").Value = 3985.3447
$ws.Range("J132").Value = 6995.7334
$ws.Range("K132").Value = 11956.0341
$ws.Range("L132").Value = 20987.2002
$ws.Range("M132").Value = -9426.034100000001
$ws.Range("N132").Value = -26047.2002
